# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the de-de / zh-cn
# handback packages have been generated:
#   - Overview "Status" column (and the per-locale "Status" column) now reads
#     "Handed back: in sync with en-US" instead of "Ready for handoff".
#   - Each locale sheet's "Latest Target File" cell now links to a.md and the
#     "Latest Handback File" / "Latest Handback DateTime" cells are populated.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: Status column (E = zh-cn, F = de-de) for both data rows.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# Columns got wider to fit the new, longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2d47d29851767a6fde02933654d918a0a7d1e468/e2e/a.md"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aMdUrl, "", "", "a.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aMdUrl, "", "", "a.md") | Out-Null

$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

$zhcn.Range("K2").Value = "2016-08-23 00:34:53"
$zhcn.Range("K3").Value = "2016-08-23 00:34:53"

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Hyperlinks.Add($dede.Range("I2"), $aMdUrl, "", "", "a.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), $aMdUrl, "", "", "a.md") | Out-Null

$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("K2").Value = "2016-08-23 00:35:04"
$dede.Range("K3").Value = "2016-08-23 00:35:04"

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
